# Update the "scraped_at" timestamps (column K) on the "snapshot" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("snapshot")

$timestamps = @{
    2  = "2025-12-14T07:01:16.200879+00:00"
    3  = "2025-12-14T07:01:16.200912+00:00"
    4  = "2025-12-14T07:01:16.200932+00:00"
    5  = "2025-12-14T07:01:18.525910+00:00"
    6  = "2025-12-14T07:01:18.525937+00:00"
    7  = "2025-12-14T07:01:18.525954+00:00"
    8  = "2025-12-14T07:01:21.323870+00:00"
    9  = "2025-12-14T07:01:24.052748+00:00"
    10 = "2025-12-14T07:01:27.062726+00:00"
    11 = "2025-12-14T07:01:29.830640+00:00"
    12 = "2025-12-14T07:01:36.908475+00:00"
    13 = "2025-12-14T07:01:36.908504+00:00"
    14 = "2025-12-14T07:01:39.219560+00:00"
    15 = "2025-12-14T07:01:41.965201+00:00"
    16 = "2025-12-14T07:01:44.307856+00:00"
    17 = "2025-12-14T07:01:46.647783+00:00"
    18 = "2025-12-14T07:01:46.647814+00:00"
    19 = "2025-12-14T07:01:46.647833+00:00"
    20 = "2025-12-14T07:01:46.647850+00:00"
    21 = "2025-12-14T07:01:46.647866+00:00"
    22 = "2025-12-14T07:01:49.524196+00:00"
    23 = "2025-12-14T07:01:49.524228+00:00"
    24 = "2025-12-14T07:01:49.524246+00:00"
    25 = "2025-12-14T07:01:52.305516+00:00"
    26 = "2025-12-14T07:01:52.305544+00:00"
    27 = "2025-12-14T07:01:52.305561+00:00"
    28 = "2025-12-14T07:01:52.305578+00:00"
    29 = "2025-12-14T07:01:52.305593+00:00"
    30 = "2025-12-14T07:01:54.652248+00:00"
    31 = "2025-12-14T07:01:54.652278+00:00"
    32 = "2025-12-14T07:01:54.652296+00:00"
    33 = "2025-12-14T07:01:56.806918+00:00"
    34 = "2025-12-14T07:01:56.806952+00:00"
    35 = "2025-12-14T07:01:56.806972+00:00"
    36 = "2025-12-14T07:01:59.547947+00:00"
    37 = "2025-12-14T07:02:02.399837+00:00"
    38 = "2025-12-14T07:02:02.399868+00:00"
    39 = "2025-12-14T07:02:02.399888+00:00"
    40 = "2025-12-14T07:02:05.260763+00:00"
    41 = "2025-12-14T07:02:07.967884+00:00"
    42 = "2025-12-14T07:02:07.967912+00:00"
    43 = "2025-12-14T07:02:10.256854+00:00"
    44 = "2025-12-14T07:02:10.256886+00:00"
}

foreach ($row in $timestamps.Keys) {
    $ws.Cells.Item($row, 11).Value = $timestamps[$row]
}
